# Auto-generated Excel COM-interop script
# Applies numeric updates to the Leve profit-calculation columns (H-N)
# across multiple worksheets, per the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 825.6
$ws.Range("I2").Value = 1136.7142
$ws.Range("J2").Value = 99.666664
$ws.Range("K2").Value = 1136.7142
$ws.Range("L2").Value = 99.666664
$ws.Range("M2").Value = -1023.7142
$ws.Range("N2").Value = -325.666664

$ws.Range("H64").Value = 4059.85
$ws.Range("I64").Value = 3844.2778
$ws.Range("K64").Value = 3844.2778
$ws.Range("M64").Value = -3596.2778

$ws.Range("H67").Value = 4059.85
$ws.Range("I67").Value = 3844.2778
$ws.Range("K67").Value = 3844.2778
$ws.Range("M67").Value = -2986.2778

$ws.Range("H129").Value = 1291.4117
$ws.Range("I129").Value = 458.15384
$ws.Range("K129").Value = 1374.46152
$ws.Range("M129").Value = 3625.53848

$ws.Range("H132").Value = 1570.075
$ws.Range("I132").Value = 1070.8529
$ws.Range("K132").Value = 3212.5587
$ws.Range("M132").Value = -682.5587000000005

$ws.Range("H138").Value = 2232.6333
$ws.Range("I138").Value = 1146.174
$ws.Range("J138").Value = 2908
$ws.Range("K138").Value = 3438.522
$ws.Range("L138").Value = 8724
$ws.Range("M138").Value = 1701.478
$ws.Range("N138").Value = -19004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2968
$ws.Range("I2").Value = 2526.5334
$ws.Range("K2").Value = 2526.5334
$ws.Range("M2").Value = -2413.5334

$ws.Range("H32").Value = 36847.2
$ws.Range("I32").Value = 21954.375
$ws.Range("K32").Value = 21954.375
$ws.Range("M32").Value = -21667.375

$ws.Range("H61").Value = 1319.6666
$ws.Range("I61").Value = 1287.7878
$ws.Range("J61").Value = 1670.3334
$ws.Range("K61").Value = 1287.7878
$ws.Range("L61").Value = 1670.3334
$ws.Range("M61").Value = -1075.7878
$ws.Range("N61").Value = -2094.3334

$ws.Range("H116").Value = 2968
$ws.Range("I116").Value = 2526.5334
$ws.Range("K116").Value = 2526.5334
$ws.Range("M116").Value = -232.5333999999998

$ws.Range("H128").Value = 77268.25
$ws.Range("J128").Value = 77268.25
$ws.Range("L128").Value = 77268.25
$ws.Range("N128").Value = -87228.25

$ws.Range("H129").Value = 100780
$ws.Range("J129").Value = 100780
$ws.Range("L129").Value = 100780
$ws.Range("N129").Value = -110780

$ws.Range("H132").Value = 8905
$ws.Range("I132").Value = 8605.259
$ws.Range("J132").Value = 16998
$ws.Range("K132").Value = 25815.777
$ws.Range("L132").Value = 50994
$ws.Range("M132").Value = -23285.777
$ws.Range("N132").Value = -56054

$ws.Range("H136").Value = 1319.6666
$ws.Range("I136").Value = 1287.7878
$ws.Range("J136").Value = 1670.3334
$ws.Range("K136").Value = 3863.3634
$ws.Range("L136").Value = 5011.0002
$ws.Range("M136").Value = -1313.3634
$ws.Range("N136").Value = -10111.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2968
$ws.Range("I3").Value = 2526.5334
$ws.Range("K3").Value = 2526.5334
$ws.Range("M3").Value = -2412.5334

$ws.Range("H43").Value = 169999
$ws.Range("J43").Value = 169999
$ws.Range("L43").Value = 169999
$ws.Range("N43").Value = -170361

$ws.Range("H94").Value = 596.6070999999999
$ws.Range("I94").Value = 612.96295
$ws.Range("J94").Value = 155
$ws.Range("K94").Value = 612.96295
$ws.Range("L94").Value = 155
$ws.Range("M94").Value = -161.96295
$ws.Range("N94").Value = -1057

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1587.4286
$ws.Range("I16").Value = 1575
$ws.Range("J16").Value = 1604
$ws.Range("K16").Value = 1575
$ws.Range("L16").Value = 1604
$ws.Range("M16").Value = -1288
$ws.Range("N16").Value = -2178

$ws.Range("H41").Value = 13999.091
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 13999.091
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 13999.091
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -14855.091

$ws.Range("H50").Value = 7499.4165
$ws.Range("J50").Value = 7499.4165
$ws.Range("L50").Value = 7499.4165
$ws.Range("N50").Value = -8749.416499999999

$ws.Range("H51").Value = 14999.077
$ws.Range("J51").Value = 14999.077
$ws.Range("L51").Value = 14999.077
$ws.Range("N51").Value = -16471.077

$ws.Range("H59").Value = 16108.667
$ws.Range("J59").Value = 16108.667
$ws.Range("L59").Value = 16108.667
$ws.Range("N59").Value = -18398.667

$ws.Range("H60").Value = 14066.134
$ws.Range("I60").Value = 1000
$ws.Range("J60").Value = 14999.429
$ws.Range("K60").Value = 1000
$ws.Range("L60").Value = 14999.429
$ws.Range("M60").Value = -489
$ws.Range("N60").Value = -16021.429

$ws.Range("H61").Value = 14999.077
$ws.Range("J61").Value = 14999.077
$ws.Range("L61").Value = 14999.077
$ws.Range("N61").Value = -15695.077

$ws.Range("H68").Value = 24990
$ws.Range("J68").Value = 24990
$ws.Range("L68").Value = 24990
$ws.Range("N68").Value = -26488

$ws.Range("H71").Value = 24990
$ws.Range("J71").Value = 24990
$ws.Range("L71").Value = 74970
$ws.Range("N71").Value = -82458

$ws.Range("H74").Value = 42060.43
$ws.Range("J74").Value = 51245
$ws.Range("L74").Value = 51245
$ws.Range("N74").Value = -52993

$ws.Range("H77").Value = 42060.43
$ws.Range("J77").Value = 51245
$ws.Range("L77").Value = 153735
$ws.Range("N77").Value = -162471

$ws.Range("H100").Value = 59000
$ws.Range("I100").Value = 58000
$ws.Range("J100").Value = 60000
$ws.Range("K100").Value = 58000
$ws.Range("L100").Value = 60000
$ws.Range("M100").Value = -56918
$ws.Range("N100").Value = -62164

$ws.Range("H113").Value = 1587.4286
$ws.Range("I113").Value = 1575
$ws.Range("J113").Value = 1604
$ws.Range("K113").Value = 1575
$ws.Range("L113").Value = 1604
$ws.Range("M113").Value = 595
$ws.Range("N113").Value = -5944

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 449999.5
$ws.Range("I128").Value = 449999.5
$ws.Range("K128").Value = 1349998.5
$ws.Range("M128").Value = -1345018.5

$ws.Range("H132").Value = 2523.3713
$ws.Range("J132").Value = 2878.125
$ws.Range("L132").Value = 25903.125
$ws.Range("N132").Value = -30963.125

$ws.Range("H137").Value = 6670928
$ws.Range("J137").Value = 6470.75
$ws.Range("L137").Value = 19412.25
$ws.Range("N137").Value = -29612.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3284.1667
$ws.Range("I122").Value = 3284.1667
$ws.Range("K122").Value = 9852.500100000001
$ws.Range("M122").Value = -7402.500100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H135").Value = 127647.06
$ws.Range("I135").Value = 50000
$ws.Range("J135").Value = 314000
$ws.Range("K135").Value = 50000
$ws.Range("L135").Value = 314000
$ws.Range("M135").Value = -44930
$ws.Range("N135").Value = -324140

$ws.Range("H136").Value = 4475.0415
$ws.Range("I136").Value = 3672.6428
$ws.Range("J136").Value = 5598.4
$ws.Range("K136").Value = 11017.9284
$ws.Range("L136").Value = 16795.2
$ws.Range("M136").Value = -8467.928400000001
$ws.Range("N136").Value = -21895.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 19490.3
$ws.Range("J41").Value = 20695.143
$ws.Range("L41").Value = 20695.143
$ws.Range("N41").Value = -21475.143

$ws.Range("H107").Value = 26318086
$ws.Range("I107").Value = 2203.3845
$ws.Range("J107").Value = 83335830
$ws.Range("K107").Value = 6610.1535
$ws.Range("L107").Value = 250007490
$ws.Range("M107").Value = -4690.1535
$ws.Range("N107").Value = -250011330

$ws.Range("H113").Value = 1229.0714
$ws.Range("I113").Value = 584.6667
$ws.Range("J113").Value = 1712.375
$ws.Range("K113").Value = 1754.0001
$ws.Range("L113").Value = 5137.125
$ws.Range("M113").Value = 415.9999
$ws.Range("N113").Value = -9477.125

$ws.Range("H126").Value = 9308
$ws.Range("I126").Value = 3030.5
$ws.Range("J126").Value = 13493
$ws.Range("K126").Value = 9091.5
$ws.Range("L126").Value = 40479
$ws.Range("M126").Value = -6621.5
$ws.Range("N126").Value = -45419
